# (SVBF-162) - QA - After code review feedback - depend on SVBF-133 - QA - After code review feedback
#
# Renames the three worksheets to prefix them with their tab order, and
# moves the "active"/selected sheet + selection around to match the new
# reviewed state of the workbook:
#   - Profiles          -> 0_Profiles        (becomes the active/selected sheet)
#   - Parameters         -> 1_Parameters      (loses the "active" tab state)
#   - UserManagement      -> 2_UserManagement  (selection moves to D23 / view
#                                               scrolls right toward column J)

$wb = $excel.ActiveWorkbook

$wsProfiles   = $wb.Worksheets.Item(1)
$wsParameters = $wb.Worksheets.Item(2)
$wsUsers      = $wb.Worksheets.Item(3)

$wsProfiles.Name   = "0_Profiles"
$wsParameters.Name = "1_Parameters"
$wsUsers.Name      = "2_UserManagement"

# Update the selection/scroll position on the UserManagement sheet before
# leaving it, so it is recorded as D23 (with the view scrolled over toward
# column J) instead of the old B3 selection.
$wsUsers.Activate() | Out-Null
$wsUsers.Range("D23").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1

# Profiles (the first sheet) becomes the active/selected tab, taking that
# state away from Parameters.
$wsProfiles.Activate() | Out-Null
